$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B73: convert it to a real number (4) instead of text "4"
$ws.Range("B73").Value = 4

# Add new row 74 with the new annotation data
$ws.Range("A74").Value = "Ying Tang"
$ws.Range("B74").Value = "'1"
$ws.Range("C74").Value = "You do not understand"
$ws.Range("D74").Value = "CRT"
$ws.Range("E74").Value = "THE"
$ws.Range("F74").Value = "9f35a425-2bea-4e69-9731-af889a0691d3"
$ws.Range("G74").Value = "r1Kr3TyAb_annotated.xlsx"
$ws.Range("H74").Value = "You do not understand the work by Veit et al."
